# Update "PUMA" sheet (sheet3 in the package, 3rd tab in the workbook):
# insert a new column before column AO (41) and populate it with the
# 2010-2014 overdose mortality rate data + headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PUMA")

# Insert a new column before the current column AO (41st column). This
# shifts the existing AO:AS data right by one, to AP:AT.
$ws.Columns.Item(41).Insert()

# Row 1 (long-form header) / Row 2 (short machine-readable header) for the
# newly inserted column.
$ws.Cells.Item(1, 41).Value = "Unintentional overdose deaths involving any drug among NYC residents age 15-84 per 100,000 people in the years 2010-2014  to all persons"
$ws.Cells.Item(2, 41).Value = "overdose_mortality_per100000_10_14"

# Data rows 3-57: 2010-2014 unintentional overdose mortality rate (all persons).
$values = @(
    7.1421000000000001,
    6.3003999999999998,
    9.4998000000000005,
    9.4245000000000001,
    16.691600000000001,
    11.1823,
    9.6609999999999996,
    11.204000000000001,
    9.0975000000000001,
    17.1264,
    7.8346999999999998,
    6.9358000000000004,
    16.707999999999998,
    15.3544,
    5.1947999999999999,
    4.6740000000000004,
    7.5799000000000003,
    4.9225000000000003,
    10.283799999999999,
    5.2385999999999999,
    19.974,
    18.209099999999999,
    15.126899999999999,
    10.9465,
    11.054399999999999,
    14.3065,
    6.9512999999999998,
    7.8324999999999996,
    9.8009000000000004,
    20.173100000000002,
    9.2443000000000008,
    6.75,
    3.9531999999999998,
    4.7359999999999998,
    6.9833999999999996,
    11.797000000000001,
    3.0386000000000002,
    5.0082000000000004,
    7.9465000000000003,
    6.9444999999999997,
    14.154299999999999,
    7.0617999999999999,
    3.7389000000000001,
    5.4659000000000004,
    7.6612999999999998,
    2.9944000000000002,
    5.7493999999999996,
    3.145,
    6.3726000000000003,
    4.1067999999999998,
    10.465299999999999,
    5.2919,
    6.4835000000000003,
    7.86,
    11.613200000000001
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(3 + $i, 41).Value = $values[$i]
}

# Restore the sheet's view: scrolled to row 28, selection on B61.
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Range("B61").Select()
